# Update gh-pages to output generated at 456a3b4
# Applies value updates + a new row insert (2024-06-01 AP动漫嘉年华 entry)
# across the 展览(1) / 演出(2) / 本地生活(3) / 全部类型(4) sheets.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$value)
    # Force literal text so date-looking strings ("2024-06-01") are not
    # auto-converted to Excel date serials, then drop the temporary
    # NumberFormat/quote-prefix styling so the cell keeps the sheet's
    # default (unstyled) look.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 8131
$ws1.Range("F3").Value = 125
$ws1.Range("F5").Value = 31344
$ws1.Range("G5").Value = 68
$ws1.Range("F8").Value = 705
$ws1.Range("F10").Value = 145
$ws1.Range("F11").Value = 448
$ws1.Range("F12").Value = 806
$ws1.Range("F13").Value = 56
$ws1.Range("F14").Value = 614
$ws1.Range("F15").Value = 401
$ws1.Range("F16").Value = 25
$ws1.Range("F17").Value = 561
$ws1.Range("F19").Value = 417
$ws1.Range("F20").Value = 420
$ws1.Range("F21").Value = 1115
$ws1.Range("F23").Value = 719
$ws1.Range("F24").Value = 2357
$ws1.Range("F25").Value = 840
$ws1.Range("F26").Value = 70
$ws1.Range("F27").Value = 1093
$ws1.Range("F29").Value = 638

# Insert new row 30 (shifts old row 30 -> row 31), copy column-A's
# bordered/bold style from the row above so A30 matches the sheet's
# running-index styling.
$ws1.Rows.Item(30).Insert()
$ws1.Range("A29").Copy()
$ws1.Range("A30").PasteSpecial(-4122)
$ws1.Range("A30").Value = 29
Set-TextValue $ws1.Range("B30") "2024-06-01"
$ws1.Range("C30").Value = "广州·第五届AP动漫嘉年华"
$ws1.Range("D30").Value = "西环路1号 广州岭南会展中心"
$ws1.Range("E30").Value = "2024.06.01 10:00-06.01 17:00"
$ws1.Range("F30").Value = 2
$ws1.Range("G30").Value = 55
$ws1.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$ws1.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

# Old row 30 (622 排球少年only) is now row 31; its running index (col A)
# advances to match the new row number, and its F value changed too.
$ws1.Range("A31").Value = 30
$ws1.Range("F31").Value = 1080

# ---------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 68

# ---------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 528

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (combined view of sheets 1-3, sorted by date)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 528
$ws4.Range("F3").Value = 8131
$ws4.Range("F4").Value = 125
$ws4.Range("F7").Value = 31347
$ws4.Range("G7").Value = 68
$ws4.Range("F10").Value = 705
$ws4.Range("F12").Value = 68
$ws4.Range("F13").Value = 145
$ws4.Range("F14").Value = 448
$ws4.Range("F18").Value = 806
$ws4.Range("F19").Value = 56
$ws4.Range("F20").Value = 614
$ws4.Range("F21").Value = 401
$ws4.Range("F23").Value = 25
$ws4.Range("F27").Value = 561
$ws4.Range("F29").Value = 417
$ws4.Range("F30").Value = 420
$ws4.Range("F31").Value = 1115
$ws4.Range("F33").Value = 719
$ws4.Range("F34").Value = 2357
$ws4.Range("F35").Value = 840
$ws4.Range("F36").Value = 70
$ws4.Range("F37").Value = 1093
$ws4.Range("F40").Value = 638

# Insert new row 41 (shifts old row 41 -> row 42).
$ws4.Rows.Item(41).Insert()
$ws4.Range("A40").Copy()
$ws4.Range("A41").PasteSpecial(-4122)
$ws4.Range("A41").Value = 40
Set-TextValue $ws4.Range("B41") "2024-06-01"
$ws4.Range("C41").Value = "广州·第五届AP动漫嘉年华"
$ws4.Range("D41").Value = "西环路1号 广州岭南会展中心"
$ws4.Range("E41").Value = "2024.06.01 10:00-06.01 17:00"
$ws4.Range("F41").Value = 2
$ws4.Range("G41").Value = 55
$ws4.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$ws4.Range("I41").Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

# Old row 41 (622 排球少年only) is now row 42; its running index (col A)
# advances to match the new row number, and its F value changed too.
$ws4.Range("A42").Value = 41
$ws4.Range("F42").Value = 1080
